# "Re-rendered sprites for better contrast"
# Tracking Table.xlsx edits:
#  - mark M102 as Done ("x")
#  - fill in missing Life/HP data for row 104 (Cost/RunningCost recalc automatically)
#  - insert a new "FC Hibberd Planet" vehicle row before the "Drewry Shunter" row,
#    pushing "Drewry Shunter" down to row 113
#  - tidy up the conditional formatting on column C to cover the newly inserted row
#  - update the saved view (selected cell / scroll position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Mark M102 as done ---
$ws.Range("M102").Value = "x"

# --- Row 104: add Life (I) / HP (J) values and correct Top Speed (G); K104/L104 recalc ---
$ws.Range("G104").Value = 62
$ws.Range("I104").Value = 25
$ws.Range("J104").Value = 375
$ws.Range("M104").Value = "x"

# --- Insert a new row above the "Drewry Shunter" row (row 112), shifting it to row 113 ---
$ws.Rows("112:112").Insert(-4121)  # xlShiftDown

# Recreate the C column "Gap to Previous" formula across both rows 112:113 in one
# shot so Excel builds it as a shared formula (same pattern used for the rest of
# the column) and each row correctly references the row directly above it.
$ws.Range("C112:C113").Formula = "=B112-B111"

# --- Populate the new row 112 with the "FC Hibberd Planet" data ---
$ws.Range("A112").Value = "FC Hibberd Planet"
$ws.Range("B112").Value = 1948
$ws.Range("D112").Value = 1
$ws.Range("E112").Value = "Cargo Tram"
$ws.Range("F112").Value = 16046
$ws.Range("G112").Value = 25
$ws.Range("H112").Value = 60
$ws.Range("I112").Value = 30
$ws.Range("J112").Value = 123
$ws.Range("K112").Formula = "=SQRT(G112*H112)*POWER((MIN(I112,20)+SQRT(MAX(I112-20,0))),0.9)*`$B`$1"
$ws.Range("L112").Formula = "=POWER((G112*G112*H112), 0.33)*LOG10(J112)*10*`$B`$1"
$ws.Range("M112").Value = "x"
$ws.Range("N112").Value = 6
$ws.Range("Q112").Formula = "=CONCATENATE(ROUND(N112*VLOOKUP(E112,'ID Scheme'!`$A`$2:`$E`$7,3),0), ""x"",ROUND(O112*VLOOKUP(E112,'ID Scheme'!`$A`$2:`$E`$7,5),0), ""x"",ROUND(P112*VLOOKUP(E112,'ID Scheme'!`$A`$2:`$E`$7,4),0))"

# --- Update the conditional formatting on column C to account for the inserted row ---
# Drop the two old single-row rules (one for the old C111, one for the old C112)...
$ws.Range("C112").FormatConditions.Item(1).Delete()
$oldC111Rule = $ws.Range("C111").FormatConditions.Item(1)
$ruleFontColor = $oldC111Rule.Font.Color
$ruleFillColor = $oldC111Rule.Interior.Color
$oldC111Rule.Delete()

# ...replace with a single rule spanning the merged C111:C113 block...
$mergedRule = $ws.Range("C111:C113").FormatConditions.Add(1, 5, "10")
$mergedRule.Font.Color = $ruleFontColor
$mergedRule.Interior.Color = $ruleFillColor

# ...and shift the tail of the big "rest of the column" rule from C113 down to C114
# (the head C1:C2 / C4:C110 portions are untouched).
$mainRule = $ws.Range("C1").FormatConditions.Item(1)
$mainFontColor = $mainRule.Font.Color
$mainFillColor = $mainRule.Interior.Color
$mainRule.ModifyAppliesToRange($ws.Range("C1:C2"))
$midRule = $ws.Range("C4:C110").FormatConditions.Add(1, 5, "10")
$midRule.Font.Color = $mainFontColor
$midRule.Interior.Color = $mainFillColor
$tailRule = $ws.Range("C114:C1048576").FormatConditions.Add(1, 5, "10")
$tailRule.Font.Color = $mainFontColor
$tailRule.Interior.Color = $mainFillColor

# --- Update the saved view: scroll position / selected cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 94
$ws.Range("J110").Select()
